# Insert a new weekly record as row 28 in the "Alcachofa" sheet.
# This pushes the former rows 28-44 down to rows 29-45 (dimension grows
# from A1:R44 to A1:R45) and populates the newly opened row 28 with a
# fresh Española / Provincia de Limarí entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 28..44 down one position, opening up a blank row 28.
$ws.Rows.Item(28).Insert()

# Fill in the new row 28 with the latest weekly price record.
$ws.Range("A28").Value = 8
$ws.Range("B28").Value = "Terminal La Palmera de La Serena"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 45205
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 100112013
$ws.Range("G28").Value = "Alcachofa"
$ws.Range("H28").Value = "Española"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 520
$ws.Range("K28").Value = 8000
$ws.Range("L28").Value = 9000
$ws.Range("M28").Value = 8500
$ws.Range("N28").Value = "`$/caja 30 unidades"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 283
$ws.Range("Q28").Value = 30
$ws.Range("R28").Value = "Hortaliza"
